$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 42604.890821759262
$ws.Range("B6").Value = "Named"
$ws.Range("C6").Value = 6902
$ws.Range("D6").Value = 4229
$ws.Range("E6").Value = 318
$ws.Range("F6").Value = 43
$ws.Range("G6").Value = 23
$ws.Range("H6").Value = 64
$ws.Range("I6").Value = 34
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 100
